$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 1156.375
$ws.Range("I12").Value = 717.1667
$ws.Range("K12").Value = 717.1667
$ws.Range("M12").Value = -547.1667

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 4374.9434
$ws.Range("J17").Value = 4374.9434
$ws.Range("L17").Value = 13124.8302
$ws.Range("N17").Value = -13460.8302

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 944.4286
$ws.Range("I58").Value = 522.2
$ws.Range("J58").Value = 2000
$ws.Range("K58").Value = 1566.6
$ws.Range("L58").Value = 6000
$ws.Range("M58").Value = -1416.6
$ws.Range("N58").Value = -6300

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 448.7
$ws.Range("I103").Value = 355
$ws.Range("J103").Value = 542.4
$ws.Range("K103").Value = 1065
$ws.Range("L103").Value = 1627.2
$ws.Range("M103").Value = -479
$ws.Range("N103").Value = -2799.2

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 2386.2917
$ws.Range("J112").Value = 2442.5
$ws.Range("L112").Value = 7327.5
$ws.Range("N112").Value = -9543.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 4234.5
$ws.Range("I113").Value = 3238.625
$ws.Range("J113").Value = 5230.375
$ws.Range("K113").Value = 3238.625
$ws.Range("L113").Value = 5230.375
$ws.Range("M113").Value = 15.375
$ws.Range("N113").Value = -11738.375

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 9262304
$ws.Range("I137").Value = 3605.25
$ws.Range("J137").Value = 27779704
$ws.Range("K137").Value = 10815.75
$ws.Range("L137").Value = 83339112
$ws.Range("M137").Value = -8265.75
$ws.Range("N137").Value = -83344212

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2731.2622
$ws.Range("I138").Value = 2409.8572
$ws.Range("J138").Value = 2900
$ws.Range("K138").Value = 7229.571599999999
$ws.Range("L138").Value = 8700
$ws.Range("M138").Value = -2089.571599999999
$ws.Range("N138").Value = -18980

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H43").Value = 43447.332
$ws.Range("I43").Value = 40342
$ws.Range("J43").Value = 45000
$ws.Range("K43").Value = 40342
$ws.Range("L43").Value = 45000
$ws.Range("M43").Value = -40029
$ws.Range("N43").Value = -45626

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 4364.7744
$ws.Range("I122").Value = 3832.0527
$ws.Range("J122").Value = 5208.25
$ws.Range("K122").Value = 11496.1581
$ws.Range("L122").Value = 15624.75
$ws.Range("M122").Value = -9046.158100000001
$ws.Range("N122").Value = -20524.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3917.3333
$ws.Range("I132").Value = 4014.1875
$ws.Range("J132").Value = 3142.5
$ws.Range("K132").Value = 12042.5625
$ws.Range("L132").Value = 9427.5
$ws.Range("M132").Value = -9512.5625
$ws.Range("N132").Value = -14487.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H138").Value = 79989.5
$ws.Range("I138").Value = 60000
$ws.Range("J138").Value = 99979
$ws.Range("K138").Value = 60000
$ws.Range("L138").Value = 99979
$ws.Range("M138").Value = -54860
$ws.Range("N138").Value = -110259

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5561412
$ws.Range("I31").Value = 12822066
$ws.Range("K31").Value = 12822066
$ws.Range("M31").Value = -12821771

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 5561412
$ws.Range("I34").Value = 12822066
$ws.Range("K34").Value = 12822066
$ws.Range("M34").Value = -12821864

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 3264.4443
$ws.Range("I99").Value = 3095
$ws.Range("K99").Value = 3095
$ws.Range("M99").Value = -1597

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 2128.3635
$ws.Range("I105").Value = 1815.6428
$ws.Range("K105").Value = 1815.6428
$ws.Range("M105").Value = -68.64280000000008

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 3264.4443
$ws.Range("I126").Value = 3095
$ws.Range("K126").Value = 9285
$ws.Range("M126").Value = -6815

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 999
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 625056
$ws.Range("I12").Value = 9.800000000000001
$ws.Range("J12").Value = 909167.9399999999
$ws.Range("K12").Value = 29.4
$ws.Range("L12").Value = 2727503.82
$ws.Range("M12").Value = 143.6
$ws.Range("N12").Value = -2727849.82

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 116.375
$ws.Range("J23").Value = 85.85714
$ws.Range("L23").Value = 257.57142
$ws.Range("N23").Value = -727.57142

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 1911.4
$ws.Range("J98").Value = 500
$ws.Range("L98").Value = 1500
$ws.Range("N98").Value = -4496

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 3128
$ws.Range("I113").Value = 2574.6667
$ws.Range("J113").Value = 3460
$ws.Range("K113").Value = 7724.000100000001
$ws.Range("L113").Value = 10380
$ws.Range("M113").Value = -5554.000100000001
$ws.Range("N113").Value = -14720

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 1359.75
$ws.Range("I114").Value = 2173.25
$ws.Range("J114").Value = 546.25
$ws.Range("K114").Value = 6519.75
$ws.Range("L114").Value = 1638.75
$ws.Range("M114").Value = -3265.75
$ws.Range("N114").Value = -8146.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1637.0769
$ws.Range("J131").Value = 1667.7727
$ws.Range("L131").Value = 5003.3181
$ws.Range("N131").Value = -15083.3181

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 999
$ws.Range("I135").Value = 0
$ws.Range("K135").Value = 0
$ws.Range("M135").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2104.4546
$ws.Range("I122").Value = 1884.2858
$ws.Range("K122").Value = 5652.857400000001
$ws.Range("M122").Value = -3202.857400000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2698.1538
$ws.Range("I132").Value = 2698.1538
$ws.Range("K132").Value = 8094.4614
$ws.Range("M132").Value = -5564.4614

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H133").Value = 95000
$ws.Range("J133").Value = 95000
$ws.Range("L133").Value = 95000
$ws.Range("N133").Value = -105120

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3705.04
$ws.Range("I40").Value = 3443.8572
$ws.Range("K40").Value = 3443.8572
$ws.Range("M40").Value = -3307.8572

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1967.8182
$ws.Range("I68").Value = 2143.375
$ws.Range("J68").Value = 1499.6666
$ws.Range("K68").Value = 2143.375
$ws.Range("L68").Value = 1499.6666
$ws.Range("M68").Value = -1394.375
$ws.Range("N68").Value = -2997.6666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 1967.8182
$ws.Range("I71").Value = 2143.375
$ws.Range("J71").Value = 1499.6666
$ws.Range("K71").Value = 10716.875
$ws.Range("L71").Value = 7498.333000000001
$ws.Range("M71").Value = -6972.875
$ws.Range("N71").Value = -14986.333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H127").Value = 25000
$ws.Range("J127").Value = 25000
$ws.Range("L127").Value = 25000
$ws.Range("N127").Value = -34920

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 19274.092
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H101").Value = 92360.60000000001
$ws.Range("J101").Value = 92360.60000000001
$ws.Range("L101").Value = 92360.60000000001
$ws.Range("N101").Value = -98850.60000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H125").Value = 59997.5
$ws.Range("J125").Value = 59997.5
$ws.Range("L125").Value = 59997.5
$ws.Range("N125").Value = -69837.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2606.8572
$ws.Range("I136").Value = 2083.0625
$ws.Range("K136").Value = 6249.1875
$ws.Range("M136").Value = -3699.1875
